# Applies scheduled-runner market price/profit updates to Jenova_Profits workbook.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets targeted cell updates
# for currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) on specific rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2134.204
$ws.Range("J17").Value = 2134.204
$ws.Range("L17").Value = 6402.612000000001
$ws.Range("N17").Value = -6738.612000000001
$ws.Range("H33").Value = 200
$ws.Range("I33").Value = 200
$ws.Range("K33").Value = 200
$ws.Range("M33").Value = 29
$ws.Range("H135").Value = 1251688.5
$ws.Range("I135").Value = 1668334.6
$ws.Range("K135").Value = 15015011.4
$ws.Range("M135").Value = -15012476.4
$ws.Range("H138").Value = 6858.278
$ws.Range("J138").Value = 7128.9375
$ws.Range("L138").Value = 21386.8125
$ws.Range("N138").Value = -31666.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 594.0909
$ws.Range("I2").Value = 603.8946999999999
$ws.Range("J2").Value = 532
$ws.Range("K2").Value = 603.8946999999999
$ws.Range("L2").Value = 532
$ws.Range("M2").Value = -490.8946999999999
$ws.Range("N2").Value = -758
$ws.Range("H32").Value = 3807.5479
$ws.Range("I32").Value = 3626.0845
$ws.Range("K32").Value = 3626.0845
$ws.Range("M32").Value = -3339.0845
$ws.Range("H45").Value = 2426.5454
$ws.Range("I45").Value = 1698.4
$ws.Range("K45").Value = 1698.4
$ws.Range("M45").Value = -1321.4
$ws.Range("H61").Value = 4517.9287
$ws.Range("I61").Value = 4080.6365
$ws.Range("K61").Value = 4080.6365
$ws.Range("M61").Value = -3868.6365
$ws.Range("H74").Value = 1702.5
$ws.Range("I74").Value = 1036.125
$ws.Range("J74").Value = 3035.25
$ws.Range("K74").Value = 1036.125
$ws.Range("L74").Value = 3035.25
$ws.Range("M74").Value = -162.125
$ws.Range("N74").Value = -4783.25
$ws.Range("H77").Value = 1702.5
$ws.Range("I77").Value = 1036.125
$ws.Range("J77").Value = 3035.25
$ws.Range("K77").Value = 5180.625
$ws.Range("L77").Value = 15176.25
$ws.Range("M77").Value = -812.625
$ws.Range("N77").Value = -23912.25
$ws.Range("H88").Value = 1934.8182
$ws.Range("I88").Value = 4212.4287
$ws.Range("J88").Value = 871.93335
$ws.Range("K88").Value = 4212.4287
$ws.Range("L88").Value = 871.93335
$ws.Range("M88").Value = -3806.4287
$ws.Range("N88").Value = -1683.93335
$ws.Range("H91").Value = 1934.8182
$ws.Range("I91").Value = 4212.4287
$ws.Range("J91").Value = 871.93335
$ws.Range("K91").Value = 4212.4287
$ws.Range("L91").Value = 871.93335
$ws.Range("M91").Value = -2808.4287
$ws.Range("N91").Value = -3679.93335
$ws.Range("H116").Value = 594.0909
$ws.Range("I116").Value = 603.8946999999999
$ws.Range("J116").Value = 532
$ws.Range("K116").Value = 603.8946999999999
$ws.Range("L116").Value = 532
$ws.Range("M116").Value = 1690.1053
$ws.Range("N116").Value = -5120
$ws.Range("H136").Value = 4517.9287
$ws.Range("I136").Value = 4080.6365
$ws.Range("K136").Value = 12241.9095
$ws.Range("M136").Value = -9691.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 594.0909
$ws.Range("I3").Value = 603.8946999999999
$ws.Range("J3").Value = 532
$ws.Range("K3").Value = 603.8946999999999
$ws.Range("L3").Value = 532
$ws.Range("M3").Value = -489.8946999999999
$ws.Range("N3").Value = -760
$ws.Range("H86").Value = 1136252.8
$ws.Range("I86").Value = 1548227.5
$ws.Range("J86").Value = 3322.25
$ws.Range("K86").Value = 1548227.5
$ws.Range("L86").Value = 3322.25
$ws.Range("M86").Value = -1547104.5
$ws.Range("N86").Value = -5568.25
$ws.Range("H89").Value = 1136252.8
$ws.Range("I89").Value = 1548227.5
$ws.Range("J89").Value = 3322.25
$ws.Range("K89").Value = 7741137.5
$ws.Range("L89").Value = 16611.25
$ws.Range("M89").Value = -7735521.5
$ws.Range("N89").Value = -27843.25
$ws.Range("H108").Value = 99996.336
$ws.Range("J108").Value = 99996.336
$ws.Range("L108").Value = 99996.336
$ws.Range("N108").Value = -107676.336
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8673.5
$ws.Range("I16").Value = 4064.6667
$ws.Range("K16").Value = 4064.6667
$ws.Range("M16").Value = -3777.6667
$ws.Range("H22").Value = 353
$ws.Range("J22").Value = 650
$ws.Range("L22").Value = 650
$ws.Range("N22").Value = -1350
$ws.Range("H113").Value = 8673.5
$ws.Range("I113").Value = 4064.6667
$ws.Range("K113").Value = 4064.6667
$ws.Range("M113").Value = -1894.6667
$ws.Range("H125").Value = 98000
$ws.Range("J125").Value = 98000
$ws.Range("L125").Value = 98000
$ws.Range("N125").Value = -102920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 39934.89
$ws.Range("J107").Value = 81815
$ws.Range("L107").Value = 245445
$ws.Range("N107").Value = -249285
$ws.Range("H121").Value = 668426.5600000001
$ws.Range("J121").Value = 1430790
$ws.Range("L121").Value = 4292370
$ws.Range("N121").Value = -4294990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 744255.8
$ws.Range("J80").Value = 772190.3
$ws.Range("L80").Value = 772190.3
$ws.Range("N80").Value = -774186.3
$ws.Range("H83").Value = 744255.8
$ws.Range("J83").Value = 772190.3
$ws.Range("L83").Value = 3860951.5
$ws.Range("N83").Value = -3870935.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3013.5557
$ws.Range("I82").Value = 3075.4
$ws.Range("J82").Value = 2936.25
$ws.Range("K82").Value = 3075.4
$ws.Range("L82").Value = 2936.25
$ws.Range("M82").Value = -2714.4
$ws.Range("N82").Value = -3658.25
$ws.Range("H85").Value = 3013.5557
$ws.Range("I85").Value = 3075.4
$ws.Range("J85").Value = 2936.25
$ws.Range("K85").Value = 3075.4
$ws.Range("L85").Value = 2936.25
$ws.Range("M85").Value = -1827.4
$ws.Range("N85").Value = -5432.25
$ws.Range("H132").Value = 10164.667
$ws.Range("I132").Value = 7874.75
$ws.Range("J132").Value = 11996.6
$ws.Range("K132").Value = 23624.25
$ws.Range("L132").Value = 35989.8
$ws.Range("M132").Value = -21094.25
$ws.Range("N132").Value = -41049.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 254471.5
$ws.Range("J96").Value = 4693.5
$ws.Range("L96").Value = 4693.5
$ws.Range("N96").Value = -7439.5
$ws.Range("H126").Value = 2134.6667
$ws.Range("I126").Value = 2100
$ws.Range("K126").Value = 6300
$ws.Range("M126").Value = -3830
